# The "eyes" / corneal-ulcer rule rows (rows 36-39) are being removed from
# the Conditions list. Clear their contents and reset the row height back
# to the sheet default (the rows had an explicit ht="30" from their old
# wrapped text; AutoFit drops that back to the default height now that the
# cells are empty), then leave the now-empty block selected, matching the
# post-edit workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A36:D39")
$rng.ClearContents()
$rng.EntireRow.AutoFit()
$rng.Select()
